$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Replace the static "entrance_datetime" values (rows 2-21) with a live =TODAY() formula
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Formula = "=TODAY()"
}

# 2) Add the new 16th student on row 22 (copy row 21's formatting first so the
#    new row matches the look of the existing data rows)
$ws.Range("A21:E21").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value = "STD000016"
$ws.Range("B22").Value = "rere"
$ws.Range("C22").Value = "Smith"
$ws.Range("D22").Formula = "=CONCATENATE(""user"",A22,""@hei.school"")"
$ws.Range("E22").Formula = "=TODAY()"

# 3) Tweak the date format used by the entrance_datetime column
$ws.Range("E2:E22").NumberFormat = "dd/mm/yy"

# 4) Leave the selection where the editor last left it
$ws.Range("F22").Select()
